$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) on both sheets to remove the test-case specific suffix
$newProductName = "4281-MS-EI-DB-SAR-REC-CTRFD-RNI-FEE-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update short name (B2) on the input sheet from a numeric literal to a text value
$wsInput.Range("B2").Value = "428q"

# Move the selection on the input sheet from B17 to B3
$wsInput.Range("B3").Select()

# Make ProductLoanOutput the active (selected) sheet/tab instead of ProductLoanInput
$wsOutput.Activate()
